$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00863423648946369
$ws.Range("C2").Value = 1.081231980802479
$ws.Range("D2").Value = 8.112271609790723
$ws.Range("E2").Value = 2.848204980297367
$ws.Range("F2").Value = 2.87966458212119
$ws.Range("G2").Value = 46

$ws.Range("B3").Value = 0.1079425770541491
$ws.Range("C3").Value = 1.025681875173506
$ws.Range("D3").Value = 5.546319209064285
$ws.Range("E3").Value = 2.355062463941092
$ws.Range("F3").Value = 2.379171181174347
$ws.Range("G3").Value = 45

$ws.Range("B4").Value = 0.08673112263663806
$ws.Range("C4").Value = 0.8989375388443066
$ws.Range("D4").Value = 4.47623325838249
$ws.Range("E4").Value = 2.115711052668225
$ws.Range("F4").Value = 2.138371923145345
$ws.Range("G4").Value = 44

$ws.Range("B5").Value = 0.1560044769001684
$ws.Range("C5").Value = 0.9224880001651093
$ws.Range("D5").Value = 4.879102887750935
$ws.Range("E5").Value = 2.20886914228773
$ws.Range("F5").Value = 2.229429331798652
$ws.Range("G5").Value = 43

$ws.Range("B6").Value = 0.1600298176880285
$ws.Range("C6").Value = 0.9424896015931804
$ws.Range("D6").Value = 4.916621931739378
$ws.Range("E6").Value = 2.217345695136277
$ws.Range("F6").Value = 2.238371143106592
$ws.Range("G6").Value = 42

$ws.Range("B7").Value = 0.1728584914046674
$ws.Range("C7").Value = 0.9206488348686618
$ws.Range("D7").Value = 4.635874527679912
$ws.Range("E7").Value = 2.153108108683796
$ws.Range("F7").Value = 2.172819442882898
$ws.Range("G7").Value = 41

$ws.Range("B8").Value = 0.1929609596835913
$ws.Range("C8").Value = 0.9736223895218249
$ws.Range("D8").Value = 4.804027380007986
$ws.Range("E8").Value = 2.191809156840072
$ws.Range("F8").Value = 2.211112598008703
$ws.Range("G8").Value = 40

$ws.Range("B9").Value = 0.1960041159357992
$ws.Range("C9").Value = 0.9531248790116683
$ws.Range("D9").Value = 4.870747427708349
$ws.Range("E9").Value = 2.206976988486366
$ws.Range("F9").Value = 2.226992678098346
$ws.Range("G9").Value = 39

$ws.Range("B10").Value = 0.2155726075138289
$ws.Range("C10").Value = 0.9906079526445191
$ws.Range("D10").Value = 5.010930836497054
$ws.Range("E10").Value = 2.238510852441206
$ws.Range("F10").Value = 2.258015469991629
$ws.Range("G10").Value = 38

$ws.Range("B11").Value = 0.2079488633544275
$ws.Range("C11").Value = 0.9496862135200643
$ws.Range("D11").Value = 5.054960332149459
$ws.Range("E11").Value = 2.248323893959555
$ws.Range("F11").Value = 2.269566474070086
$ws.Range("G11").Value = 37

